# código refactorizado de ejercicio apache poi
#
# - Hoja1!B1 now holds the "7 results have been found." label (reuses the
#   existing shared string), extending row 1's used range to column B.
# - The previously-selected sheet (Hoja1) loses focus; "Credenciales"
#   becomes the active/selected sheet in the workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("B1").Value = "7 results have been found."

$wsCred = $wb.Worksheets.Item("Credenciales")
$wsCred.Activate()
